$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("2021-04-02", "01:38:03", 0.4350838195794973),
    @("2021-04-02", "01:39:21", 0.4339909234046336),
    @("2021-04-02", "01:40:05", 0.4339909234046336),
    @("2021-04-02", "01:41:07", 0.4328980272297699)
)

$row = 11
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = "'" + $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $row++
}
